$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New data rows to append (Nombre, edad, genero, residencia, sisben)
$data = @(
    @("Valentina", 20, "F", "Rural",  2),
    @("Sandra",    21, "F", "Urbano", 7),
    @("Yuly",      22, "F", "Urbano", 1),
    @("Yuri",      23, "F", "Urbano", 2),
    @("Daniela",   24, "F", "Urbano", 10)
)

$startRow = 16
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $values = $data[$i]

    # Copy formatting from the row above so borders/fills match existing data rows
    $ws.Range("A" + ($row - 1) + ":E" + ($row - 1)).Copy() | Out-Null
    $ws.Range("A" + $row + ":E" + $row).PasteSpecial(-4122) | Out-Null

    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
    $ws.Cells.Item($row, 4).Value = $values[3]
    $ws.Cells.Item($row, 5).Value = $values[4]
}

$excel.CutCopyMode = 0

# Update selection / view to match recorded state
$ws.Range("F23").Select() | Out-Null
$ws.Application.ActiveWindow.ScrollRow = 4

$wb.Save()
